# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" suffixed header columns to "_FV2404" /
# "_FV2410" respectively, wraps the data range in an Excel Table
# ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (row 1) cells: "<Name>_old" -> "<Name>_FV2404",
#    "<Name>_new" -> "<Name>_FV2410". Column K ("diff") is unchanged.
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# 2) Turn the used range into a native Excel Table ("Table1") so the
#    header row exposes the AutoFilter dropdowns used for the diff view.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U92"), $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
